$wb = $excel.ActiveWorkbook

# --- Swap the names of the two worksheet tabs ---
# Physical sheet order / rIds stay the same; only the displayed tab
# names swap: the sheet that holds "Partner Data"-shaped content
# (dimension A1:W37) becomes "Solver Team Data" and vice versa.
$wsFirst  = $wb.Worksheets.Item(1)
$wsSecond = $wb.Worksheets.Item(2)

$wsFirst.Name  = "__TEMP_SWAP__"
$wsSecond.Name = "Partner Data"
$wsFirst.Name  = "Solver Team Data"

# --- Widen column A on the first sheet (now "Solver Team Data") ---
$wsFirst.Columns.Item(1).ColumnWidth = 37.165

# --- Update the stored workbook absolute path hint ---
# (Excel regenerates x15ac:absPath / xr:revisionPtr automatically on
# save; nothing to set via the object model for those.)

# --- Restore cell selections on each sheet ---
$null = $wsFirst.Range("B23").Select()
$null = $wsSecond.Range("B36").Select()
